$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 19, shifting existing rows 19-121 down to 20-122
$ws.Rows.Item(19).Insert()

# Fill in the new row 19 with the latest week's data
$ws.Cells.Item(19, 1).Value = 11
$ws.Cells.Item(19, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(19, 3).Value = "Bíobío"
$ws.Cells.Item(19, 4).Value = 44547
$ws.Cells.Item(19, 4).NumberFormat = $ws.Cells.Item(20, 4).NumberFormat
$ws.Cells.Item(19, 5).Value = 8
$ws.Cells.Item(19, 6).Value = 100112003
$ws.Cells.Item(19, 7).Value = "Ajo"
$ws.Cells.Item(19, 8).Value = "Chino"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 200
$ws.Cells.Item(19, 11).Value = 17000
$ws.Cells.Item(19, 12).Value = 18000
$ws.Cells.Item(19, 13).Value = 17500
$ws.Cells.Item(19, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(19, 15).Value = "China"
$ws.Cells.Item(19, 16).Value = 1750
$ws.Cells.Item(19, 17).Value = 10
$ws.Cells.Item(19, 18).Value = "Hortaliza"
